# Update schematics and BOM
# Replace the generic "Thick Film Resistors - SMD 0603" description with the
# specific part descriptions for each resistor BOM line (R1, R2/R3, R4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe keeps these as literal text (preserves the existing
# quote-prefixed cell style instead of Excel minting a new one).
$ws.Range("D5").Formula = "'Res Thick Film 0603 100 Ohm 1% 0.25W(1/4W) ±100ppm/C Pad SMD Automotive T/R"
$ws.Range("D6").Formula = "'Res Thick Film 0603 3.9K Ohm 1% 1/10W ±100ppm/°C Molded SMD SMD Paper T/R"
$ws.Range("D7").Formula = "'Res Thick Film 0603 1K Ohm 1% 1/10W ±100ppm/°C Molded SMD Punched Carrier T/R"

# Move the active selection to H8 (matches the saved cursor position).
[void]$ws.Range("H8").Select()
